# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets, per the latest scrape.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1099
$ws1.Range("F4").Value = 1753
$ws1.Range("F5").Value = 775
$ws1.Range("F6").Value = 96

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1099
$ws4.Range("F4").Value = 1753
$ws4.Range("F6").Value = 775
$ws4.Range("F7").Value = 96
